$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style_D2 = $ws.Range("D2").Style
$ws.Range("D2").Value = "'58.139.05"
$ws.Range("D2").Style = $style_D2
$style_E2 = $ws.Range("E2").Style
$ws.Range("E2").Value = "'  -4.16%  "
$ws.Range("E2").Style = $style_E2

$style_D3 = $ws.Range("D3").Style
$ws.Range("D3").Value = "'2.725.69"
$ws.Range("D3").Style = $style_D3
$style_E3 = $ws.Range("E3").Style
$ws.Range("E3").Value = "'  -6.41%  "
$ws.Range("E3").Style = $style_E3

$style_E4 = $ws.Range("E4").Style
$ws.Range("E4").Value = "'  +0.06%  "
$ws.Range("E4").Style = $style_E4

$style_D5 = $ws.Range("D5").Style
$ws.Range("D5").Value = "'493.96"
$ws.Range("D5").Style = $style_D5
$style_E5 = $ws.Range("E5").Style
$ws.Range("E5").Value = "'  -6.55%  "
$ws.Range("E5").Style = $style_E5

$style_D6 = $ws.Range("D6").Style
$ws.Range("D6").Value = "'132.94"
$ws.Range("D6").Style = $style_D6
$style_E6 = $ws.Range("E6").Style
$ws.Range("E6").Value = "'  -6.86%  "
$ws.Range("E6").Style = $style_E6

$style_E7 = $ws.Range("E7").Style
$ws.Range("E7").Value = "'  +0.10%  "
$ws.Range("E7").Style = $style_E7

$style_D8 = $ws.Range("D8").Style
$ws.Range("D8").Value = "'0.522"
$ws.Range("D8").Style = $style_D8
$style_E8 = $ws.Range("E8").Style
$ws.Range("E8").Value = "'  -5.30%  "
$ws.Range("E8").Style = $style_E8

$style_D9 = $ws.Range("D9").Style
$ws.Range("D9").Value = "'2.734.14"
$ws.Range("D9").Style = $style_D9
$style_E9 = $ws.Range("E9").Style
$ws.Range("E9").Value = "'  -6.11%  "
$ws.Range("E9").Style = $style_E9

$style_D10 = $ws.Range("D10").Style
$ws.Range("D10").Value = "'5.84"
$ws.Range("D10").Style = $style_D10
$style_E10 = $ws.Range("E10").Style
$ws.Range("E10").Value = "'  -0.15%  "
$ws.Range("E10").Style = $style_E10

$style_E11 = $ws.Range("E11").Style
$ws.Range("E11").Value = "'  -6.22%  "
$ws.Range("E11").Style = $style_E11

$style_D12 = $ws.Range("D12").Style
$ws.Range("D12").Value = "'0.340"
$ws.Range("D12").Style = $style_D12
$style_E12 = $ws.Range("E12").Style
$ws.Range("E12").Value = "'  -3.50%  "
$ws.Range("E12").Style = $style_E12

$style_D13 = $ws.Range("D13").Style
$ws.Range("D13").Value = "'0.126"
$ws.Range("D13").Style = $style_D13
$style_E13 = $ws.Range("E13").Style
$ws.Range("E13").Value = "'  +1.19%  "
$ws.Range("E13").Style = $style_E13

$style_D14 = $ws.Range("D14").Style
$ws.Range("D14").Value = "'3.219.74"
$ws.Range("D14").Style = $style_D14
$style_E14 = $ws.Range("E14").Style
$ws.Range("E14").Value = "'  -5.79%  "
$ws.Range("E14").Style = $style_E14

$style_D15 = $ws.Range("D15").Style
$ws.Range("D15").Value = "'58.257.20"
$ws.Range("D15").Style = $style_D15
$style_E15 = $ws.Range("E15").Style
$ws.Range("E15").Value = "'  -4.14%  "
$ws.Range("E15").Style = $style_E15

$style_D16 = $ws.Range("D16").Style
$ws.Range("D16").Value = "'20.98"
$ws.Range("D16").Style = $style_D16
$style_E16 = $ws.Range("E16").Style
$ws.Range("E16").Value = "'  -7.28%  "
$ws.Range("E16").Style = $style_E16

$style_D17 = $ws.Range("D17").Style
$ws.Range("D17").Value = "'2.738.15"
$ws.Range("D17").Style = $style_D17
$style_E17 = $ws.Range("E17").Style
$ws.Range("E17").Value = "'  -5.76%  "
$ws.Range("E17").Style = $style_E17

$style_E18 = $ws.Range("E18").Style
$ws.Range("E18").Value = "'  -6.00%  "
$ws.Range("E18").Style = $style_E18

$style_D19 = $ws.Range("D19").Style
$ws.Range("D19").Value = "'4.62"
$ws.Range("D19").Style = $style_D19
$style_E19 = $ws.Range("E19").Style
$ws.Range("E19").Value = "'  -6.31%  "
$ws.Range("E19").Style = $style_E19

$style_D20 = $ws.Range("D20").Style
$ws.Range("D20").Value = "'340.87"
$ws.Range("D20").Style = $style_D20
$style_E20 = $ws.Range("E20").Style
$ws.Range("E20").Value = "'  -5.35%  "
$ws.Range("E20").Style = $style_E20

$style_D21 = $ws.Range("D21").Style
$ws.Range("D21").Value = "'10.71"
$ws.Range("D21").Style = $style_D21
$style_E21 = $ws.Range("E21").Style
$ws.Range("E21").Value = "'  -7.07%  "
$ws.Range("E21").Style = $style_E21

$style_D22 = $ws.Range("D22").Style
$ws.Range("D22").Value = "'6.09"
$ws.Range("D22").Style = $style_D22
$style_E22 = $ws.Range("E22").Style
$ws.Range("E22").Value = "'  -5.67%  "
$ws.Range("E22").Style = $style_E22

$style_D23 = $ws.Range("D23").Style
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = $style_D23
$style_E23 = $ws.Range("E23").Style
$ws.Range("E23").Value = "'  -0.08%  "
$ws.Range("E23").Style = $style_E23

$style_E24 = $ws.Range("E24").Style
$ws.Range("E24").Value = "'  -0.79%  "
$ws.Range("E24").Style = $style_E24

$style_D25 = $ws.Range("D25").Style
$ws.Range("D25").Value = "'61.74"
$ws.Range("D25").Style = $style_D25
$style_E25 = $ws.Range("E25").Style
$ws.Range("E25").Value = "'  -2.51%  "
$ws.Range("E25").Style = $style_E25

$style_D26 = $ws.Range("D26").Style
$ws.Range("D26").Value = "'0.418"
$ws.Range("D26").Style = $style_D26
$style_E26 = $ws.Range("E26").Style
$ws.Range("E26").Value = "'  -6.93%  "
$ws.Range("E26").Style = $style_E26

$style_D27 = $ws.Range("D27").Style
$ws.Range("D27").Value = "'0.168"
$ws.Range("D27").Style = $style_D27
$style_E27 = $ws.Range("E27").Style
$ws.Range("E27").Value = "'  -7.58%  "
$ws.Range("E27").Style = $style_E27

$style_E28 = $ws.Range("E28").Style
$ws.Range("E28").Value = "'  +0.14%  "
$ws.Range("E28").Style = $style_E28

$style_D29 = $ws.Range("D29").Style
$ws.Range("D29").Value = "'7.19"
$ws.Range("D29").Style = $style_D29
$style_E29 = $ws.Range("E29").Style
$ws.Range("E29").Value = "'  -5.96%  "
$ws.Range("E29").Style = $style_E29

$style_D30 = $ws.Range("D30").Style
$ws.Range("D30").Value = "'0.0₃0785"
$ws.Range("D30").Style = $style_D30
$style_E30 = $ws.Range("E30").Style
$ws.Range("E30").Value = "'  -8.34%  "
$ws.Range("E30").Style = $style_E30

$style_E31 = $ws.Range("E31").Style
$ws.Range("E31").Value = "'  -0.03%  "
$ws.Range("E31").Style = $style_E31

$style_D32 = $ws.Range("D32").Style
$ws.Range("D32").Value = "'1.57"
$ws.Range("D32").Style = $style_D32
$style_E32 = $ws.Range("E32").Style
$ws.Range("E32").Value = "'  -5.64%  "
$ws.Range("E32").Style = $style_E32

$style_D33 = $ws.Range("D33").Style
$ws.Range("D33").Value = "'18.73"
$ws.Range("D33").Style = $style_D33
$style_E33 = $ws.Range("E33").Style
$ws.Range("E33").Value = "'  -4.55%  "
$ws.Range("E33").Style = $style_E33

$style_D34 = $ws.Range("D34").Style
$ws.Range("D34").Value = "'147.61"
$ws.Range("D34").Style = $style_D34
$style_E34 = $ws.Range("E34").Style
$ws.Range("E34").Value = "'  -3.86%  "
$ws.Range("E34").Style = $style_E34

$style_D35 = $ws.Range("D35").Style
$ws.Range("D35").Value = "'4.04"
$ws.Range("D35").Style = $style_D35
$style_E35 = $ws.Range("E35").Style
$ws.Range("E35").Value = "'  -6.74%  "
$ws.Range("E35").Style = $style_E35

$style_D36 = $ws.Range("D36").Style
$ws.Range("D36").Value = "'5.22"
$ws.Range("D36").Style = $style_D36
$style_E36 = $ws.Range("E36").Style
$ws.Range("E36").Value = "'  -6.31%  "
$ws.Range("E36").Style = $style_E36

$style_D37 = $ws.Range("D37").Style
$ws.Range("D37").Value = "'0.891"
$ws.Range("D37").Style = $style_D37
$style_E37 = $ws.Range("E37").Style
$ws.Range("E37").Value = "'  -10.86%  "
$ws.Range("E37").Style = $style_E37

$style_D38 = $ws.Range("D38").Style
$ws.Range("D38").Value = "'1.10"
$ws.Range("D38").Style = $style_D38
$style_E38 = $ws.Range("E38").Style
$ws.Range("E38").Value = "'  -8.26%  "
$ws.Range("E38").Style = $style_E38

$style_D39 = $ws.Range("D39").Style
$ws.Range("D39").Value = "'36.34"
$ws.Range("D39").Style = $style_D39
$style_E39 = $ws.Range("E39").Style
$ws.Range("E39").Value = "'  -4.21%  "
$ws.Range("E39").Style = $style_E39

$style_D40 = $ws.Range("D40").Style
$ws.Range("D40").Value = "'2.148.06"
$ws.Range("D40").Style = $style_D40
$style_E40 = $ws.Range("E40").Style
$ws.Range("E40").Value = "'  -8.01%  "
$ws.Range("E40").Style = $style_E40

$style_D41 = $ws.Range("D41").Style
$ws.Range("D41").Value = "'0.999"
$ws.Range("D41").Style = $style_D41
$style_E41 = $ws.Range("E41").Style
$ws.Range("E41").Value = "'  +0.10%  "
$ws.Range("E41").Style = $style_E41

$style_D42 = $ws.Range("D42").Style
$ws.Range("D42").Value = "'3.44"
$ws.Range("D42").Style = $style_D42
$style_E42 = $ws.Range("E42").Style
$ws.Range("E42").Value = "'  -6.35%  "
$ws.Range("E42").Style = $style_E42

$style_D43 = $ws.Range("D43").Style
$ws.Range("D43").Value = "'0.0546"
$ws.Range("D43").Style = $style_D43
$style_E43 = $ws.Range("E43").Style
$ws.Range("E43").Value = "'  -3.69%  "
$ws.Range("E43").Style = $style_E43

$style_E44 = $ws.Range("E44").Style
$ws.Range("E44").Value = "'  -9.96%  "
$ws.Range("E44").Style = $style_E44

$style_D45 = $ws.Range("D45").Style
$ws.Range("D45").Value = "'0.585"
$ws.Range("D45").Style = $style_D45
$style_E45 = $ws.Range("E45").Style
$ws.Range("E45").Value = "'  -8.94%  "
$ws.Range("E45").Style = $style_E45

$style_E46 = $ws.Range("E46").Style
$ws.Range("E46").Value = "'  -0.14%  "
$ws.Range("E46").Style = $style_E46

$style_D47 = $ws.Range("D47").Style
$ws.Range("D47").Value = "'18.59"
$ws.Range("D47").Style = $style_D47
$style_E47 = $ws.Range("E47").Style
$ws.Range("E47").Value = "'  -10.47%  "
$ws.Range("E47").Style = $style_E47

$style_E48 = $ws.Range("E48").Style
$ws.Range("E48").Value = "'  -5.25%  "
$ws.Range("E48").Style = $style_E48

$style_D49 = $ws.Range("D49").Style
$ws.Range("D49").Value = "'4.47"
$ws.Range("D49").Style = $style_D49
$style_E49 = $ws.Range("E49").Style
$ws.Range("E49").Value = "'  -7.24%  "
$ws.Range("E49").Style = $style_E49

$style_D50 = $ws.Range("D50").Style
$ws.Range("D50").Value = "'0.0874"
$ws.Range("D50").Style = $style_D50
$style_E50 = $ws.Range("E50").Style
$ws.Range("E50").Value = "'  -5.06%  "
$ws.Range("E50").Style = $style_E50

$style_D51 = $ws.Range("D51").Style
$ws.Range("D51").Value = "'17.01"
$ws.Range("D51").Style = $style_D51
$style_E51 = $ws.Range("E51").Style
$ws.Range("E51").Value = "'  -7.28%  "
$ws.Range("E51").Style = $style_E51
